$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 155, shifting existing rows 155-235 down to 156-236.
$ws.Rows("155:155").Insert()

# Populate the newly inserted row 155 with the new weekly data point.
$ws.Range("A155").Value = 6
$ws.Range("B155").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C155").Value = "Metropolitana"
$ws.Range("D155").Value = 44813
$ws.Range("E155").Value = 13
$ws.Range("F155").Value = 100112029
$ws.Range("G155").Value = "Orégano"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 52
$ws.Range("K155").Value = 15000
$ws.Range("L155").Value = 16000
$ws.Range("M155").Value = 15442
$ws.Range("N155").Value = "$/docena de atados"
$ws.Range("O155").Value = "Región Metropolitana"
$ws.Range("P155").Value = 5147
$ws.Range("Q155").Value = 3
$ws.Range("R155").Value = "Hortaliza"

$ws.Range("D155").NumberFormat = "YYYY-MM-DD HH:MM:SS"
